$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.704230555368213
$ws.Range("C2").Value = 0.7115922907102572
$ws.Range("D2").Value = 0.06462893496600586
$ws.Range("E2").Value = 0.126717272406907
$ws.Range("F2").Value = 2.986501631082632
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 1.363216454538374
$ws.Range("J2").Value = 0.2204158715211832
$ws.Range("B3").Value = 1.586601651986484
$ws.Range("C3").Value = 0.6620414870534432
$ws.Range("D3").Value = 0.06429010423659065
$ws.Range("E3").Value = 0.1246604846344717
$ws.Range("F3").Value = 2.938505611647486
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 1.350534565481667
$ws.Range("J3").Value = 0.215419335374591
$ws.Range("B4").Value = 1.515334518319605
$ws.Range("C4").Value = 0.6320526742333641
$ws.Range("D4").Value = 0.06409995989444539
$ws.Range("E4").Value = 0.1234666764559051
$ws.Range("F4").Value = 2.910970408700209
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 1.343614453065243
$ws.Range("J4").Value = 0.2124921238356663
$ws.Range("B5").Value = 1.486532062299716
$ws.Range("C5").Value = 0.6199403482800676
$ws.Range("D5").Value = 0.06402702871518784
$ws.Range("E5").Value = 0.1229975236956768
$ws.Range("F5").Value = 2.900232230864759
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 1.341010528148999
$ws.Range("J5").Value = 0.2113344033914828
$ws.Range("B6").Value = 1.481763864135303
$ws.Range("C6").Value = 0.6179356167092465
$ws.Range("D6").Value = 0.06401519519098642
$ws.Range("E6").Value = 0.122920667204145
$ws.Range("F6").Value = 2.898478203933479
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 1.340591145240332
$ws.Range("J6").Value = 0.2111442805295667
$ws.Range("B7").Value = 1.514945110331269
$ws.Range("C7").Value = 0.6318888859561298
$ws.Range("D7").Value = 0.0640989578083726
$ws.Range("E7").Value = 0.1234602791568271
$ws.Range("F7").Value = 2.910823640692897
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 1.343578463156391
$ws.Range("J7").Value = 0.212476368429428
$ws.Range("B8").Value = 1.663472547840229
$ws.Range("C8").Value = 0.6944160393873631
$ws.Range("D8").Value = 0.06450842268690948
$ws.Range("E8").Value = 0.125993735023922
$ws.Range("F8").Value = 2.969548761537482
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 1.358662755243699
$ws.Range("J8").Value = 0.2186637237227274
$ws.Range("B9").Value = 1.962409578618633
$ws.Range("C9").Value = 0.8205485658610314
$ws.Range("D9").Value = 0.06545117601579165
$ws.Range("E9").Value = 0.1315118951701102
$ws.Range("F9").Value = 3.100240432245783
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 1.395205310009516
$ws.Range("J9").Value = 0.2319248999536256
$ws.Range("B10").Value = 2.186854196268257
$ws.Range("C10").Value = 0.9154579630124999
$ws.Range("D10").Value = 0.06622621222041403
$ws.Range("E10").Value = 0.1359049820292739
$ws.Range("F10").Value = 3.205994922406006
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 1.426421110286142
$ws.Range("J10").Value = 0.2423728667186822
$ws.Range("B11").Value = 2.290037890962651
$ws.Range("C11").Value = 0.9591429206125781
$ws.Range("D11").Value = 0.06659611178709213
$ws.Range("E11").Value = 0.1379779485213604
$ws.Range("F11").Value = 3.256278357772828
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 1.441597586395034
$ws.Range("J11").Value = 0.2472829179435792
$ws.Range("B12").Value = 2.329268775899891
$ws.Range("C12").Value = 0.9757602543107282
$ws.Range("D12").Value = 0.06673862651983598
$ws.Range("E12").Value = 0.1387737041176109
$ws.Range("F12").Value = 3.275636784759627
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 1.447487059101348
$ws.Range("J12").Value = 0.249165134036474
$ws.Range("B13").Value = 2.320812670831742
$ws.Range("C13").Value = 0.9721780648827689
$ws.Range("D13").Value = 0.06670782561718624
$ws.Range("E13").Value = 0.1386018441569092
$ws.Range("F13").Value = 3.271453423834743
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 1.446212285892841
$ws.Range("J13").Value = 0.2487587427600886
$ws.Range("B14").Value = 2.293262273344112
$ws.Range("C14").Value = 0.9605085314994426
$ws.Range("D14").Value = 0.06660778790845256
$ws.Range("E14").Value = 0.1380431997022384
$ws.Range("F14").Value = 3.257864605058842
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 1.442079249507671
$ws.Range("J14").Value = 0.247437308705031
$ws.Range("B15").Value = 2.276407419070324
$ws.Range("C15").Value = 0.9533703878877873
$ws.Range("D15").Value = 0.06654682843423387
$ws.Range("E15").Value = 0.1377024177326405
$ws.Range("F15").Value = 3.249582498917107
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 1.439566261556735
$ws.Range("J15").Value = 0.2466308810145676
$ws.Range("B16").Value = 2.180132835181212
$ws.Range("C16").Value = 0.9126134354957571
$ws.Range("D16").Value = 0.0662023828362166
$ws.Range("E16").Value = 0.1357710128125618
$ws.Range("F16").Value = 3.202752926076101
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 1.425449114453656
$ws.Range("J16").Value = 0.2420551706396736
$ws.Range("B17").Value = 2.121350035276237
$ws.Range("C17").Value = 0.8877420390747943
$ws.Range("D17").Value = 0.06599547948132312
$ws.Range("E17").Value = 0.1346052826548032
$ws.Range("F17").Value = 3.174584724031405
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 1.417040221467417
$ws.Range("J17").Value = 0.2392885843727584
$ws.Range("B18").Value = 2.087641522741421
$ws.Range("C18").Value = 0.8734846043908533
$ws.Range("D18").Value = 0.06587810844639108
$ws.Range("E18").Value = 0.1339417973469246
$ws.Range("F18").Value = 3.158587592075406
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 1.412295412065319
$ws.Range("J18").Value = 0.2377120983199887
$ws.Range("B19").Value = 2.076245831244421
$ws.Range("C19").Value = 0.8686654721012701
$ws.Range("D19").Value = 0.0658386508768487
$ws.Range("E19").Value = 0.1337183554010437
$ws.Range("F19").Value = 3.153206223576802
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 1.410704598339521
$ws.Range("J19").Value = 0.2371808575220484
$ws.Range("B20").Value = 2.127597016315121
$ws.Range("C20").Value = 0.8903846704394596
$ws.Range("D20").Value = 0.06601733592063397
$ws.Range("E20").Value = 0.1347286505596301
$ws.Range("F20").Value = 3.177562083442524
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 1.417925850269626
$ws.Range("J20").Value = 0.2395815604074159
$ws.Range("B21").Value = 2.301350205025301
$ws.Range("C21").Value = 0.9639341150633527
$ws.Range("D21").Value = 0.06663710553696944
$ws.Range("E21").Value = 0.1382069944566666
$ws.Range("F21").Value = 3.261847327021513
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 1.44328933912972
$ws.Range("J21").Value = 0.2478248225463204
$ws.Range("B22").Value = 2.415826630459208
$ws.Range("C22").Value = 1.012439540589696
$ws.Range("D22").Value = 0.06705637066156811
$ws.Range("E22").Value = 0.1405430838728066
$ws.Range("F22").Value = 3.318783462534299
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 1.460697256997435
$ws.Range("J22").Value = 0.2533458026805704
$ws.Range("B23").Value = 2.35464374359691
$ws.Range("C23").Value = 0.9865108519377372
$ws.Range("D23").Value = 0.06683131722039803
$ws.Range("E23").Value = 0.1392905061520651
$ws.Range("F23").Value = 3.288224755190413
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 1.451329554171679
$ws.Range("J23").Value = 0.2503868400755351
$ws.Range("B24").Value = 2.124772486885092
$ws.Range("C24").Value = 0.8891898077597489
$ws.Range("D24").Value = 0.06600744969792771
$ws.Range("E24").Value = 0.1346728550381755
$ws.Range("F24").Value = 3.176215405528922
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 1.417525178638954
$ws.Range("J24").Value = 0.2394490621652068
$ws.Range("B25").Value = 1.880703683865192
$ws.Range("C25").Value = 0.7860400888337722
$ws.Range("D25").Value = 0.06518146600801344
$ws.Range("E25").Value = 0.129959839104238
$ws.Range("F25").Value = 3.063194701740883
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 1.384561511340095
$ws.Range("J25").Value = 0.2282148771385124
